$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 193, shifting rows 193:239 down to 194:240
$ws.Rows.Item(193).Insert()

# Populate the new row 193 with the new record's data.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R are identical to the rest of the data set in this sheet,
# so copy them from the (now shifted down) row 194 which still holds the old row193 data,
# then set the differing columns D,J,K,L,M,P explicitly.
$ws.Range("A193:R193").Value2 = $ws.Range("A194:R194").Value2

$ws.Cells.Item(193, 4).Value = 44543
$ws.Cells.Item(193, 10).Value = 160
$ws.Cells.Item(193, 11).Value = 2500
$ws.Cells.Item(193, 12).Value = 2500
$ws.Cells.Item(193, 13).Value = 2500
$ws.Cells.Item(193, 16).Value = 833
